$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.020.73'
$ws.Range('E2').Value = '  +4.96%  '
$ws.Range('D3').Value = '2.617.52'
$ws.Range('E3').Value = '  +5.50%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.21'
$ws.Range('E5').Value = '  +3.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '181.73'
$ws.Range('E6').Value = '  +4.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').Value = '2.616.45'
$ws.Range('E9').Value = '  +5.48%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  +14.65%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.348'
$ws.Range('E12').Value = '  +4.84%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.06'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.76'
$ws.Range('E14').Value = '  +6.19%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.018.27'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000183'
$ws.Range('E16').Value = '  +8.14%  '
$ws.Range('D17').Value = '70.985.99'
$ws.Range('E17').Value = '  +5.03%  '
$ws.Range('D18').Value = '2.609.92'
$ws.Range('E18').Value = '  +5.44%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '383.16'
$ws.Range('E19').Value = '  +10.60%  '
$ws.Range('E20').Value = '  +6.98%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.49'
$ws.Range('E21').Value = '  +6.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.19'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '72.15'
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.43'
$ws.Range('E24').Value = '  +6.24%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  +11.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.72'
$ws.Range('E27').Value = '  +10.77%  '
$ws.Range('D28').Value = '2.748.52'
$ws.Range('E28').Value = '  +5.57%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  +7.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '527.91'
$ws.Range('E31').Value = '  +7.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.07'
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('E33').Value = '  +7.18%  '
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.94'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.22'
$ws.Range('E38').Value = '  +5.38%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.91'
$ws.Range('E39').Value = '  +11.20%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.94'
$ws.Range('E40').Value = '  +1.74%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.37'
$ws.Range('E41').Value = '  +6.85%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('E43').Value = '  +6.74%  '
$ws.Range('E44').Value = '  +9.64%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.08'
$ws.Range('E46').Value = '  +3.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '153.91'
$ws.Range('E47').Value = '  +4.35%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.66'
$ws.Range('E48').Value = '  +4.23%  '
$ws.Range('D49').Value = '0.0₆0274'
$ws.Range('E49').Value = '  +8.45%  '
$ws.Range('E50').Value = '  +5.15%  '
$ws.Range('E51').Value = '  +7.31%  '
